# 2021.01.08 Ki Quy 3 2020 Gen - parse column last year get ki quater A,D
#
# This script reproduces, via the Excel COM object model, the user-visible
# edits captured in the commit:
#   1. Header row 5 gets explicit 0/1/2/3 values typed into Y5:Z5:AA5:AB5
#      (previously blank placeholder cells under the "Diem" sub-headers).
#   2. Two new helper columns are filled in for every data row (8-61):
#        AH -> "Nam A Quy" (year A) and AI -> "Nam D Quy" (year D) -
#      these are parsed-out "last year" quarter values per employee record.
#   3. The window is rescrolled/zoomed and the selection left on AD5 -
#      cosmetic view state that Excel persists into the sheetView.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Row 5 "sequence" values under the score sub-headers (Y5:AB5 were blank)
# ---------------------------------------------------------------------------
$ws.Range("Y5").Value  = 0
$ws.Range("Z5").Value  = 1
$ws.Range("AA5").Value = 2
$ws.Range("AB5").Value = 3

# ---------------------------------------------------------------------------
# 2. New AH / AI columns: year parsed out of the "Ki quy" data, per row
#    (row, AH value, AI value)
# ---------------------------------------------------------------------------
$yearData = @(
    @(8, 2018, 2019),
    @(9, 2018, 2019),
    @(10, 2018, 2019),
    @(11, 2018, 2019),
    @(12, 2017, 2019),
    @(13, 2017, 2018),
    @(14, 2017, 2018),
    @(15, 2018, 2018),
    @(16, 2018, 2020),
    @(17, 2018, 2020),
    @(18, 2018, 2020),
    @(19, 2019, 2020),
    @(20, 2019, 2020),
    @(21, 2019, 2018),
    @(22, 2019, 2018),
    @(23, 2018, 2018),
    @(24, 2018, 2018),
    @(25, 2020, 2018),
    @(26, 2020, 2018),
    @(27, 2020, 2018),
    @(28, 2020, 2020),
    @(29, 2018, 2020),
    @(30, 2018, 2020),
    @(31, 2018, 2020),
    @(32, 2018, 2018),
    @(33, 2018, 2018),
    @(34, 2018, 2018),
    @(35, 2018, 2018),
    @(36, 2018, 2018),
    @(37, 2018, 2019),
    @(38, 2018, 2019),
    @(39, 2018, 2019),
    @(40, 2018, 2019),
    @(41, 2019, 2019),
    @(42, 2019, 2018),
    @(43, 2019, 2018),
    @(44, 2019, 2018),
    @(45, 2019, 2018),
    @(46, 2018, 2018),
    @(47, 2018, 2018),
    @(48, 2018, 2018),
    @(49, 2018, 2020),
    @(50, 2018, 2020),
    @(51, 2018, 2020),
    @(52, 2018, 2020),
    @(53, 2018, 2020),
    @(54, 2020, 2020),
    @(55, 2020, 2018),
    @(56, 2020, 2018),
    @(57, 2020, 2018),
    @(58, 2018, 2018),
    @(59, 2018, 2019),
    @(60, 2018, 2018),
    @(61, 2018, 2018)
)

foreach ($entry in $yearData) {
    $r  = $entry[0]
    $ah = $entry[1]
    $ai = $entry[2]
    $ws.Cells.Item($r, 34).Value = $ah   # column AH
    $ws.Cells.Item($r, 35).Value = $ai   # column AI
}

# ---------------------------------------------------------------------------
# 3. View state: zoom to 70%, rest the view near the top of the sheet and
#    leave the selection on AD5
# ---------------------------------------------------------------------------
$ws.Range("A4").Select()
$excel.ActiveWindow.Zoom = 70
$ws.Range("AD5").Select()
